$d = $word.ActiveDocument
foreach ($story in $d.StoryRanges) {
    if ($story.StoryType -eq 9 -and $story.InlineShapes.Count -gt 0) {
        $s = $story.InlineShapes.Item(1)
        Write-Output ("shape range start=" + $s.Range.Start)
        try {
            $s.Name = "image1.png"
            Write-Output "OK"
        } catch {
            Write-Output ("EXC: " + $_.Exception.Message)
        }
    }
}
